$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-06-15 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-16 Monday", 2)

$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $tbl.Cell($row, $col)
    $cell.Range.Text = $newText
}

# Row 1 (problem block 1)
Set-CellText 1 1 "340÷8="
Set-CellText 1 2 "687÷6="
Set-CellText 1 3 "178÷2="
Set-CellText 1 4 "638÷8="
Set-CellText 1 5 "832÷2="

# Row 5 (problem block 2)
Set-CellText 5 1 "764÷2="
Set-CellText 5 2 "411÷6="
Set-CellText 5 3 "742÷5="
Set-CellText 5 4 "982÷2="
Set-CellText 5 5 "235÷2="

# Row 9 (problem block 3)
Set-CellText 9 1 "710÷8="
Set-CellText 9 2 "603÷2="
Set-CellText 9 3 "756÷5="
Set-CellText 9 4 "523÷4="
Set-CellText 9 5 "986÷2="

# Row 13 (problem block 4)
Set-CellText 13 1 "301÷6="
Set-CellText 13 2 "573÷2="
Set-CellText 13 3 "944÷4="
Set-CellText 13 4 "925÷4="
Set-CellText 13 5 "882÷9="

# Row 17 (problem block 5)
Set-CellText 17 1 "994÷8="
Set-CellText 17 2 "968÷7="
Set-CellText 17 3 "298÷9="
Set-CellText 17 4 "241÷3="
Set-CellText 17 5 "282÷5="
